$d = $word.ActiveDocument

# The paragraph currently reads: <div>\n<id>p016r_a1</id>\n<head>...
# split across three runs: "<id>", "p016r_a1", "</id>".
# Collapse them into a single run reading "<id>p016r_1</id>", inheriting
# the Courier-New/gold formatting of the opening "<id>" run.
$d.Content.Find.Execute("<id>p016r_a1</id>", $false, $false, $false, $false, `
    $false, $true, 1, $false, "<id>p016r_1</id>", 2) | Out-Null
